$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: RandomForestRegressor - values updated
$ws.Range("B3").Value = 0.9823409710142269
$ws.Range("C3").Value = 0.9818270428934931
$ws.Range("D3").Value = 0.9823211752767453

# Row 4: GradientBoostingRegressor -> DecisionTreeRegressor - values updated
$ws.Range("A4").Value = "DecisionTreeRegressor"
$ws.Range("B4").Value = 0.9839127042472712
$ws.Range("C4").Value = 0.9848937939169423
$ws.Range("D4").Value = 0.9840816463475243

# Row 5: AdaBoostRegressor -> MLPRegressor - values updated
$ws.Range("A5").Value = "MLPRegressor"
$ws.Range("B5").Value = 0.8724958073020699
$ws.Range("C5").Value = 0.8687247131190465
$ws.Range("D5").Value = 0.8681428833948451
